# New crime data collected — weekly CompStat update for the 32nd Precinct.
# Updates the report header (volume/issue number + week-covering date range)
# and refreshes the crime-complaint statistics table (rows 14-30, 33) with
# the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 32   Number  5" -> "...Number  6"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  6"

# ---------------------------------------------------------------------
# Header text: reporting week date range
# "Report Covering the Week  1/27/2025  Through  2/2/2025"
#   -> "...2/3/2025  Through  2/9/2025"
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# ---------------------------------------------------------------------
# Crime-complaint statistics table: plain value updates (style unchanged)
# ---------------------------------------------------------------------
$ws.Range("N14").Value = -83.333333333333
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 500
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 700
$ws.Range("M15").Value = 166.666666666667
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -40
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -14.814814814814
$ws.Range("L16").Value = 35.294117647058
$ws.Range("M16").Value = -28.125
$ws.Range("N16").Value = -80
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -48.717948717948
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = -33.333333333333
$ws.Range("L17").Value = -26.923076923076
$ws.Range("M17").Value = 8.571428571428
$ws.Range("N17").Value = -51.898734177215
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -20
$ws.Range("M18").Value = -20
$ws.Range("N18").Value = -87.628865979381
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 133.333333333333
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 47.619047619047
$ws.Range("I19").Value = 44
$ws.Range("J19").Value = 29
$ws.Range("K19").Value = 51.724137931034
$ws.Range("L19").Value = 91.304347826087
$ws.Range("M19").Value = 41.935483870967
$ws.Range("N19").Value = -18.518518518518
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -79.411764705882
$ws.Range("C21").Value = 19
$ws.Range("E21").Value = -29.629629629629
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = -14.851485148514
$ws.Range("I21").Value = 133
$ws.Range("J21").Value = 144
$ws.Range("K21").Value = -7.638888888888
$ws.Range("L21").Value = 7.258064516129
$ws.Range("M21").Value = 9.016393442622
$ws.Range("N21").Value = -66.329113924050
$ws.Range("M22").Value = 50
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = -18.181818181818
$ws.Range("I23").Value = 25
$ws.Range("J23").Value = 31
$ws.Range("K23").Value = -19.354838709677
$ws.Range("L23").Value = 8.695652173913
$ws.Range("M23").Value = 38.888888888888
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 7.142857142857
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 13.432835820895
$ws.Range("I24").Value = 94
$ws.Range("J24").Value = 88
$ws.Range("K24").Value = 6.818181818181
$ws.Range("L24").Value = -4.081632653061
$ws.Range("M24").Value = 32.394366197183
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -23.529411764705
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 22
$ws.Range("K25").Value = -27.272727272727
$ws.Range("L25").Value = -30.434782608695
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = 2.325581395348
$ws.Range("I26").Value = 53
$ws.Range("J26").Value = 61
$ws.Range("K26").Value = -13.114754098360
$ws.Range("L26").Value = -19.696969696969
$ws.Range("M26").Value = -44.791666666666
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 166.666666666667
$ws.Range("L27").Value = 300
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -42.857142857142
$ws.Range("L28").Value = -20
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("M29").Value = -75
$ws.Range("N29").Value = -93.75
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("M30").Value = -66.666666666666
$ws.Range("N30").Value = -92.307692307692

# ---------------------------------------------------------------------
# Cells that flip between a numeric 0 (rendered via text placeholder
# shared-string "0") and a real count, or vice versa. Excel shows a
# dash/"0" via a text cell (style 13) rather than a numeric cell
# (style 14) when the count is zero, so the underlying cell type has to
# flip along with the value. We set the new value first (numbers as
# plain numbers; the text "0" placeholder via a leading apostrophe so
# it is stored as text, matching the source file's convention), then
# copy the number-format/style from an already-correct donor cell of
# the desired kind so the resulting style index lines up with the rest
# of the sheet instead of synthesizing a brand-new style.
# ---------------------------------------------------------------------

# C20: had a real count (2) -> now zero, shown as text "0" (style 13)
$ws.Range("C20").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null

# C29: had a real count (1) -> now zero, shown as text "0" (style 13)
$ws.Range("C29").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null

# C30: had a real count (1) -> now zero, shown as text "0" (style 13)
$ws.Range("C30").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null

# F33: had the text "0" placeholder (style 13) -> now a real count (1)
$ws.Range("F33").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null

# I33: had the text "0" placeholder (style 13) -> now a real count (1)
$ws.Range("I33").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("I33").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
